$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: A5 changes from 5 to 4 (B/C/D/E/F/G unchanged)
$ws.Range("A5").Value = 4

# Row 6: rebuild with new action "confirm" data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "confirm"
$ws.Range("C6").Value = "Selenium"
$ws.Range("D6").Value = "xpath"
$ws.Range("E6").Value = "//body[@id='gsr']/div[@id='main']/div[@id='cnt']/div[@class='mw']/div[@id='rcnt']/div[@class='col']/div[@id='center_col']/div[@id='res']/div[@id='search']/div/div[@id='rso']/div[1]/div[1]/div[1]/div[1]/div[1]/a[1]/h3[1]"

# Update the selection to match the diff (cursor moved to B9)
$ws.Range("B9").Select()
